# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates
$ws1.Range("F2").Value = 1164
$ws1.Range("F4").Value = 267
$ws1.Range("F5").Value = 147
$ws1.Range("F6").Value = 10
$ws1.Range("F7").Value = 12258
$ws1.Range("F8").Value = 56
$ws1.Range("F10").Value = 3
$ws1.Range("F12").Value = 12037
$ws1.Range("F13").Value = 4808
$ws1.Range("F14").Value = 2656
$ws1.Range("F15").Value = 119
$ws1.Range("F16").Value = 51
$ws1.Range("F19").Value = 946
$ws1.Range("F22").Value = 164
$ws1.Range("F23").Value = 69

# 全部类型 sheet updates (same events, combined view)
$ws4.Range("F2").Value = 1164
$ws4.Range("F4").Value = 267
$ws4.Range("F5").Value = 147
$ws4.Range("F8").Value = 10
$ws4.Range("F9").Value = 12258
$ws4.Range("F10").Value = 56
$ws4.Range("F12").Value = 3
$ws4.Range("F14").Value = 12037
$ws4.Range("F15").Value = 4808
$ws4.Range("F16").Value = 2656
$ws4.Range("F17").Value = 119
$ws4.Range("F18").Value = 51
$ws4.Range("F21").Value = 946
$ws4.Range("F24").Value = 164
$ws4.Range("F25").Value = 69

$wb.Save()
